$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 values to 2 decimal places (custom accuracy), keeping already-short values unchanged
$ws.Range("B5").Value = 23.6
$ws.Range("C5").Value = 17.43
$ws.Range("D5").Value = 0.65
$ws.Range("E5").Value = 49.5
$ws.Range("F5").Value = 41.29
$ws.Range("G5").Value = 18.43
$ws.Range("H5").Value = 70.4
$ws.Range("J5").Value = 12.53
$ws.Range("K5").Value = 19.13
$ws.Range("L5").Value = 19.98
$ws.Range("O5").Value = 17.81
$ws.Range("P5").Value = 26.08
$ws.Range("Q5").Value = 14.8
$ws.Range("S5").Value = 0.6
$ws.Range("T5").Value = 267.19
$ws.Range("U5").Value = 50.46
$ws.Range("V5").Value = 16.82
$ws.Range("W5").Value = 34.24
$ws.Range("X5").Value = 17.77
$ws.Range("Y5").Value = 2.37
$ws.Range("Z5").Value = 34.36
$ws.Range("AA5").Value = 14.69
$ws.Range("AB5").Value = 12.74
$ws.Range("AC5").Value = 15.04
$ws.Range("AD5").Value = 21.48
$ws.Range("AE5").Value = 0.52
$ws.Range("AF5").Value = 63.61
$ws.Range("AG5").Value = 9.34
$ws.Range("AH5").Value = 20.66

# Remove the last data row (row 6) entirely
$ws.Rows("6").Delete()
